# Fix #9722 - [Feature] Translate export search reports
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Produits
$ws2 = $wb.Worksheets.Item(2)   # Prélèvements
$ws3 = $wb.Worksheets.Item(3)   # Analyses

# --- Rename sheets (French -> English) ---
$ws1.Name = "Products"
$ws2.Name = "Samples"
$ws3.Name = "Analysis"

# --- Translate header row (row 3) on the "Samples" sheet ---
$ws2.Range("B3").Value = "Entity"
$ws2.Range("C3").Value = "Control point"
$ws2.Range("D3").Value = "Step"
$ws2.Range("E3").Value = "Sampler"
$ws2.Range("F3").Value = "Controller"
$ws2.Range("G3").Value = "Correcter"
$ws2.Range("H3").Value = "Sample"
$ws2.Range("I3").Value = "Date"
$ws2.Range("J3").Value = "Status"

# --- Translate header row (row 3) on the "Analysis" sheet ---
$ws3.Range("B3").Value = "Entity"
$ws3.Range("C3").Value = "Sample"
$ws3.Range("D3").Value = "Type"
$ws3.Range("E3").Value = "Temperature"
$ws3.Range("F3").Value = "Method"
$ws3.Range("G3").Value = "Caracteristics"
$ws3.Range("H3").Value = "Min"
$ws3.Range("I3").Value = "Max"
$ws3.Range("J3").Value = "Unit"
$ws3.Range("K3").Value = "Criteria"
$ws3.Range("L3").Value = "Value"
$ws3.Range("M3").Value = "Status"
$ws3.Range("N3").Value = "Release Control"

# --- Recolor header bands: green -> blue (unify with the "Products" sheet) ---
$blue = 5521920
$ws3.Range("B3:H3").Interior.Color = $blue
$ws3.Range("B3:H3").HorizontalAlignment = $ws1.Range("B3").HorizontalAlignment
$ws3.Range("I3:N3").Interior.Color = $blue
$ws2.Range("B3:J3").Interior.Color = $blue

# --- Selections on each sheet ---
$ws1.Activate()
$ws1.Range("A3").Select()

$ws2.Activate()
$ws2.Range("B3").Select()

$ws3.Activate()
$ws3.Range("A3").Select()

# "Analysis" becomes the active tab/sheet
$ws3.Activate()
